$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row=2; D=44216; M=55; N=11000; O=12000; P=11545; Q='$/caja 14 kilos empedrada'; S=825 },
    @{ Row=3; D=44181; M=65; N=9000; O=10000; P=9462; Q='$/caja 14 kilos empedrada'; S=676 },
    @{ Row=4; D=44229; M=55; N=11000; O=12000; P=11364; Q='$/caja 14 kilos empedrada'; S=812 },
    @{ Row=5; D=45236; M=80; N=18000; O=18000; P=18000; Q='$/caja 14 kilos granel'; S=1286 },
    @{ Row=6; D=45138; M=50; N=14000; O=14000; P=14000; Q='$/caja 14 kilos granel'; S=1000 },
    @{ Row=7; D=45211; M=50; N=17000; O=17000; P=17000; Q='$/caja 14 kilos granel'; S=1214 },
    @{ Row=8; D=45196; M=30; N=15000; O=15000; P=15000; Q='$/caja 14 kilos empedrada'; S=1071 },
    @{ Row=9; D=45152; M=60; N=16000; O=16000; P=16000; Q='$/caja 14 kilos empedrada'; S=1143 },
    @{ Row=10; D=45238; M=80; N=18000; O=18000; P=18000; Q='$/caja 14 kilos empedrada'; S=1286 },
    @{ Row=11; D=45140; M=30; N=15000; O=15000; P=15000; Q='$/caja 14 kilos granel'; S=1071 },
    @{ Row=12; D=45224; M=80; N=15000; O=15000; P=15000; Q='$/caja 14 kilos granel'; S=1071 },
    @{ Row=13; D=45194; M=60; N=15000; O=15000; P=15000; Q='$/caja 14 kilos granel'; S=1071 },
    @{ Row=14; D=45260; M=60; N=14000; O=14000; P=14000; Q='$/caja 14 kilos empedrada'; S=1000 },
    @{ Row=15; D=45222; M=80; N=15000; O=15000; P=15000; Q='$/caja 14 kilos granel'; S=1071 },
    @{ Row=16; D=44253; M=90; N=12000; O=13000; P=12667; Q='$/caja 14 kilos empedrada'; S=905 },
    @{ Row=17; D=45212; M=40; N=17000; O=17000; P=17000; Q='$/caja 14 kilos granel'; S=1214 },
    @{ Row=18; D=45167; M=50; N=16000; O=16000; P=16000; Q='$/caja 14 kilos empedrada'; S=1143 },
    @{ Row=19; D=44210; M=70; N=10000; O=11000; P=10357; Q='$/caja 14 kilos empedrada'; S=740 },
    @{ Row=20; D=45250; M=150; N=17000; O=17000; P=17000; Q='$/caja 14 kilos empedrada'; S=1214 },
    @{ Row=21; D=44172; M=90; N=8500; O=9000; P=8806; Q='$/caja 14 kilos empedrada'; S=629 },
    @{ Row=22; D=45240; M=50; N=16000; O=16000; P=16000; Q='$/caja 14 kilos granel'; S=1143 },
    @{ Row=23; D=44232; M=60; N=11000; O=12000; P=11583; Q='$/caja 14 kilos empedrada'; S=827 },
    @{ Row=24; D=45142; M=30; N=15000; O=15000; P=15000; Q='$/caja 14 kilos empedrada'; S=1071 },
    @{ Row=25; D=45142; M=30; N=14000; O=14000; P=14000; Q='$/caja 14 kilos granel'; S=1000 },
    @{ Row=26; D=45155; M=60; N=15000; O=15000; P=15000; Q='$/caja 14 kilos empedrada'; S=1071 }
)

foreach ($item in $rows) {
    $r = $item.Row
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 13).Value = $item.M
    $ws.Cells.Item($r, 14).Value = $item.N
    $ws.Cells.Item($r, 15).Value = $item.O
    $ws.Cells.Item($r, 16).Value = $item.P
    $ws.Cells.Item($r, 17).Value = $item.Q
    $ws.Cells.Item($r, 19).Value = $item.S
}
